# Rename sheet "Property1" -> "DataNode" (conceptually unifying
# DataNode / DataTable / Entity naming, per commit message).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Move/restore the saved cursor position to D26, matching the
# selection recorded in the sheet view on save.
$ws.Range("D26").Select()
